$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iRiaL673"
$ws.Range("B2").Value = 23110860
$ws.Range("C2").Value = "zwqhjhs73"
$ws.Range("D2").Value = 'sv9XR7$#'
$ws.Range("F2").Value = "psabtdSX"
$ws.Range("G2").Value = "TVAM"
